# Auto-generated script applying scheduled market-data updates to profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2099.9167
$ws.Range("J17").Value = 2319.9
$ws.Range("L17").Value = 6959.700000000001
$ws.Range("N17").Value = -7295.700000000001
$ws.Range("H33").Value = 177
$ws.Range("I33").Value = 181.66667
$ws.Range("K33").Value = 181.66667
$ws.Range("M33").Value = 47.33332999999999
$ws.Range("H138").Value = 3449.2058
$ws.Range("I138").Value = 4062
$ws.Range("J138").Value = 3389.9033
$ws.Range("K138").Value = 12186
$ws.Range("L138").Value = 10169.7099
$ws.Range("M138").Value = -7046
$ws.Range("N138").Value = -20449.7099
$ws.Range("H140").Value = 95999.8
$ws.Range("J140").Value = 95999.8
$ws.Range("L140").Value = 95999.8
$ws.Range("N140").Value = -106359.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 900
$ws.Range("I10").Value = 900
$ws.Range("K10").Value = 900
$ws.Range("M10").Value = -730
$ws.Range("H32").Value = 163552.86
$ws.Range("I32").Value = 223729.9
$ws.Range("J32").Value = 13110.223
$ws.Range("K32").Value = 223729.9
$ws.Range("L32").Value = 13110.223
$ws.Range("M32").Value = -223442.9
$ws.Range("N32").Value = -13684.223
$ws.Range("H122").Value = 1084.9744
$ws.Range("I122").Value = 975.1429000000001
$ws.Range("K122").Value = 2925.4287
$ws.Range("M122").Value = -475.4287000000004
$ws.Range("H132").Value = 456659.06
$ws.Range("I132").Value = 512308.03
$ws.Range("J132").Value = 2192.5
$ws.Range("K132").Value = 1536924.09
$ws.Range("L132").Value = 6577.5
$ws.Range("M132").Value = -1534394.09
$ws.Range("N132").Value = -11637.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1576.2
$ws.Range("I5").Value = 345.25
$ws.Range("K5").Value = 345.25
$ws.Range("M5").Value = -232.25
$ws.Range("H12").Value = 545.6667
$ws.Range("I12").Value = 340
$ws.Range("J12").Value = 751.3333
$ws.Range("K12").Value = 340
$ws.Range("L12").Value = 751.3333
$ws.Range("M12").Value = -172
$ws.Range("N12").Value = -1087.3333
$ws.Range("H94").Value = 1654.375
$ws.Range("I94").Value = 1664.6666
$ws.Range("K94").Value = 1664.6666
$ws.Range("M94").Value = -1213.6666
$ws.Range("H107").Value = 982.5294
$ws.Range("I107").Value = 762.6875
$ws.Range("K107").Value = 762.6875
$ws.Range("M107").Value = 1157.3125
$ws.Range("H120").Value = 4995
$ws.Range("J120").Value = 4995
$ws.Range("L120").Value = 4995
$ws.Range("N120").Value = -14671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 367.66666
$ws.Range("I5").Value = 364.83334
$ws.Range("J5").Value = 370.5
$ws.Range("K5").Value = 364.83334
$ws.Range("L5").Value = 370.5
$ws.Range("M5").Value = -252.83334
$ws.Range("N5").Value = -594.5
$ws.Range("H8").Value = 510
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H15").Value = 3999
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 3999
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 3999
$ws.Range("N15").Value = -4339
$ws.Range("M15").ClearContents()
$ws.Range("H22").Value = 565.5
$ws.Range("I22").Value = 217.85715
$ws.Range("K22").Value = 217.85715
$ws.Range("M22").Value = 132.14285
$ws.Range("H31").Value = 1503811.2
$ws.Range("I31").Value = 2139206.2
$ws.Range("K31").Value = 2139206.2
$ws.Range("M31").Value = -2138911.2
$ws.Range("H34").Value = 1503811.2
$ws.Range("I34").Value = 2139206.2
$ws.Range("K34").Value = 2139206.2
$ws.Range("M34").Value = -2139004.2
$ws.Range("H43").Value = 14751
$ws.Range("J43").Value = 14828
$ws.Range("L43").Value = 14828
$ws.Range("N43").Value = -15196
$ws.Range("H44").Value = 15001
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H58").Value = 3793132.2
$ws.Range("I58").Value = 5459.25
$ws.Range("J58").Value = 5957517
$ws.Range("K58").Value = 5459.25
$ws.Range("L58").Value = 5957517
$ws.Range("M58").Value = -5256.25
$ws.Range("N58").Value = -5957923
$ws.Range("H94").Value = 1687.9445
$ws.Range("I94").Value = 1089.2
$ws.Range("J94").Value = 1918.2307
$ws.Range("K94").Value = 1089.2
$ws.Range("L94").Value = 1918.2307
$ws.Range("M94").Value = -638.2
$ws.Range("N94").Value = -2820.2307
$ws.Range("H99").Value = 15384.588
$ws.Range("I99").Value = 22135.182
$ws.Range("J99").Value = 3008.5
$ws.Range("K99").Value = 22135.182
$ws.Range("L99").Value = 3008.5
$ws.Range("M99").Value = -20637.182
$ws.Range("N99").Value = -6004.5
$ws.Range("H101").Value = 14751
$ws.Range("J101").Value = 14828
$ws.Range("L101").Value = 14828
$ws.Range("N101").Value = -21318
$ws.Range("H126").Value = 15384.588
$ws.Range("I126").Value = 22135.182
$ws.Range("J126").Value = 3008.5
$ws.Range("K126").Value = 66405.546
$ws.Range("L126").Value = 9025.5
$ws.Range("M126").Value = -63935.546
$ws.Range("N126").Value = -13965.5
$ws.Range("H132").Value = 2681.5454
$ws.Range("I132").Value = 2609.1562
$ws.Range("K132").Value = 7827.4686
$ws.Range("M132").Value = -5297.4686
$ws.Range("H136").Value = 3793132.2
$ws.Range("I136").Value = 5459.25
$ws.Range("J136").Value = 5957517
$ws.Range("K136").Value = 16377.75
$ws.Range("L136").Value = 17872551
$ws.Range("M136").Value = -13827.75
$ws.Range("N136").Value = -17877651
$ws.Range("H141").Value = 188723.16
$ws.Range("J141").Value = 232341.2
$ws.Range("L141").Value = 232341.2
$ws.Range("N141").Value = -242701.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2075
$ws.Range("I25").Value = 2548
$ws.Range("J25").Value = 1680.8334
$ws.Range("K25").Value = 7644
$ws.Range("L25").Value = 5042.5002
$ws.Range("M25").Value = -7475
$ws.Range("N25").Value = -5380.5002
$ws.Range("H30").Value = 2075
$ws.Range("I30").Value = 2548
$ws.Range("J30").Value = 1680.8334
$ws.Range("K30").Value = 7644
$ws.Range("L30").Value = 5042.5002
$ws.Range("M30").Value = -7542
$ws.Range("N30").Value = -5246.5002
$ws.Range("H113").Value = 1256.9412
$ws.Range("J113").Value = 1430.7693
$ws.Range("L113").Value = 4292.3079
$ws.Range("N113").Value = -8632.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 9999.666999999999
$ws.Range("I12").Value = 9999.5
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 9999.5
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -9859.5
$ws.Range("N12").Value = -10280
$ws.Range("H126").Value = 2362.375
$ws.Range("I126").Value = 2362.375
$ws.Range("K126").Value = 7087.125
$ws.Range("M126").Value = -4617.125
$ws.Range("H136").Value = 65923.89
$ws.Range("J136").Value = 65923.89
$ws.Range("L136").Value = 197771.67
$ws.Range("N136").Value = -202871.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H132").Value = 5845109
$ws.Range("I132").Value = 10624418
$ws.Range("J132").Value = 3731.111
$ws.Range("K132").Value = 31873254
$ws.Range("L132").Value = 11193.333
$ws.Range("M132").Value = -31870724
$ws.Range("N132").Value = -16253.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 9500
$ws.Range("I32").Value = 9500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9500
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9183
$ws.Range("N32").ClearContents()
$ws.Range("H122").Value = 48412.832
$ws.Range("I122").Value = 1358.4706
$ws.Range("K122").Value = 4075.4118
$ws.Range("M122").Value = -1625.4118
$ws.Range("H132").Value = 7249178
$ws.Range("I132").Value = 7578490.5
$ws.Range("J132").Value = 4300
$ws.Range("K132").Value = 22735471.5
$ws.Range("L132").Value = 12900
$ws.Range("M132").Value = -22732941.5
$ws.Range("N132").Value = -17960
$ws.Range("H140").Value = 93409
$ws.Range("J140").Value = 93409
$ws.Range("L140").Value = 93409
$ws.Range("N140").Value = -103769
